$p = $ppt.ActivePresentation

$oldDate = "12/15/2018"
$newDate = "4/15/2019"

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master "datetimeFigureOut" placeholder.
$master = $p.SlideMaster
Update-DateShapes $master.Shapes

# Every slide layout off the master has its own copy of the placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# Slide 1: rename Person* shapes to Patient* (UML diagram relabel).
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $txt = $shp.TextFrame.TextRange.Text
        if ($txt -eq "PersonListPanel") {
            $shp.TextFrame.TextRange.Text = "PatientListPanel"
        } elseif ($txt -eq "PersonCard") {
            $shp.TextFrame.TextRange.Text = "PatientCard"
        }
    }
}
